# Updates cryptos list: refreshed prices / 1h volume %, and the
# Monero <-> PEPE row order swap (rows 29/30), per the Sep 3 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.912.34"
$ws.Range("E2").Value = "  +0.79%  "
# Row 3
$ws.Range("D3").Value = "2.503.17"
$ws.Range("E3").Value = "  -0.70%  "
# Row 4
$ws.Range("E4").Value = "  +0.17%  "
# Row 5
$ws.Range("D5").Value = "'534.27"
$ws.Range("E5").Value = "  +2.67%  "
# Row 6
$ws.Range("D6").Value = "'134.06"
$ws.Range("E6").Value = "  +1.63%  "
# Row 7
$ws.Range("E7").Value = "  +0.12%  "
# Row 8
$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  +2.40%  "
# Row 9
$ws.Range("D9").Value = "2.506.43"
$ws.Range("E9").Value = "  -0.51%  "
# Row 10
$ws.Range("D10").Value = "'0.0994"
$ws.Range("E10").Value = "  +2.17%  "
# Row 11
$ws.Range("E11").Value = "  -2.72%  "
# Row 12
$ws.Range("D12").Value = "'5.17"
$ws.Range("E12").Value = "  -0.76%  "
# Row 13
$ws.Range("D13").Value = "'0.330"
$ws.Range("E13").Value = "  -1.18%  "
# Row 14
$ws.Range("D14").Value = "2.946.32"
$ws.Range("E14").Value = "  -0.49%  "
# Row 15
$ws.Range("D15").Value = "58.725.53"
$ws.Range("E15").Value = "  +0.65%  "
# Row 16
$ws.Range("D16").Value = "'22.36"
$ws.Range("E16").Value = "  +0.58%  "
# Row 17
$ws.Range("E17").Value = "  +0.52%  "
# Row 18
$ws.Range("D18").Value = "2.501.62"
$ws.Range("E18").Value = "  -0.47%  "
# Row 19
$ws.Range("D19").Value = "'10.64"
$ws.Range("E19").Value = "  -0.82%  "
# Row 20
$ws.Range("D20").Value = "'4.25"
$ws.Range("E20").Value = "  +1.47%  "
# Row 21
$ws.Range("D21").Value = "'321.26"
$ws.Range("E21").Value = "  -0.87%  "
# Row 22
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +2.13%  "
# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.18%  "
# Row 24
$ws.Range("D24").Value = "'65.87"
$ws.Range("E24").Value = "  +3.23%  "
# Row 25
$ws.Range("D25").Value = "'0.409"
$ws.Range("E25").Value = "  +0.77%  "
# Row 27
$ws.Range("E27").Value = "  -1.41%  "
# Row 28
$ws.Range("D28").Value = "'7.44"
$ws.Range("E28").Value = "  +1.01%  "
# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0755"
$ws.Range("E29").Value = "  +1.04%  "
# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'172.04"
$ws.Range("E30").Value = "  +2.39%  "
# Row 31
$ws.Range("E31").Value = "  +1.58%  "
# Row 32
$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  -0.25%  "
# Row 33
$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  -0.38%  "
# Row 34
$ws.Range("E34").Value = "  +0.04%  "
# Row 35
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  -0.07%  "
# Row 36
$ws.Range("D36").Value = "'18.11"
$ws.Range("E36").Value = "  +0.36%  "
# Row 37
$ws.Range("E37").Value = "  -3.01%  "
# Row 38
$ws.Range("D38").Value = "'3.96"
$ws.Range("E38").Value = "  +0.31%  "
# Row 39
$ws.Range("E39").Value = "  +3.85%  "
# Row 40
$ws.Range("D40").Value = "'0.825"
# Row 41
$ws.Range("D41").Value = "'36.43"
$ws.Range("E41").Value = "  -1.00%  "
# Row 42
$ws.Range("D42").Value = "'3.48"
$ws.Range("E42").Value = "  +1.20%  "
# Row 43
$ws.Range("D43").Value = "'275.17"
$ws.Range("E43").Value = "  -1.33%  "
# Row 44
$ws.Range("D44").Value = "'131.17"
$ws.Range("E44").Value = "  +7.00%  "
# Row 45
$ws.Range("D45").Value = "'5.03"
$ws.Range("E45").Value = "  -0.43%  "
# Row 46
$ws.Range("E46").Value = "  -1.24%  "
# Row 47
$ws.Range("D47").Value = "'0.0935"
$ws.Range("E47").Value = "  +1.58%  "
# Row 48
$ws.Range("D48").Value = "'0.0510"
$ws.Range("E48").Value = "  +2.33%  "
# Row 49
$ws.Range("E49").Value = "  +1.78%  "
# Row 50
$ws.Range("D50").Value = "'16.82"
$ws.Range("E50").Value = "  -1.11%  "
# Row 51
$ws.Range("D51").Value = "1.749.81"
$ws.Range("E51").Value = "  +0.15%  "
